# LOM3266.xlsx update
# - "Objetivos:" row (row 10) B/C -> teacher code "519033 - Carlos Yujiro Shigue"
# - "Programa resumido:" row (row 13) B/C -> "01/01/2023" (kept as TEXT, not a date)
# - "Programa:" row (row 15) B/C -> teacher code "519033 - Carlos Yujiro Shigue"
# - "Método:" row (row 18) B/C -> teacher code "7290967 - Emerson Gonçalves de Melo"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# "01/01/2023" looks like a date to Excel, so force text entry by pre-formatting
# the cells as Text, then restore the original (General) number format/style by
# copying it over from the untouched row above once the literal text is in place.
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"
